$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "34% IE" -> "34% Internet Explorer"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("34% IE", $false, $false, $false, $false, $false, `
    $true, 1, $false, "34% Internet Explorer", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) "33% Firefox" -> three runs with identical formatting:
#       "33% " + "Mozilla " + "Firefox"
#    A plain Find/Replace would just rewrite the single run's text, so
#    instead we locate the run, split off a point after "33% " and
#    insert "Mozilla " there. Toggling Bold on the freshly-inserted
#    text (on, then back off) forces the engine to keep it as its own
#    run instead of silently re-merging it with the neighbouring runs
#    that share the same resulting formatting.
# ---------------------------------------------------------------------
$match = $d.Content
$found = $match.Find.Execute("33% Firefox", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPoint = $match.Duplicate
    $splitPoint.Collapse(1)
    $splitPoint.MoveEnd(1, 4)     # "33% " is 4 characters
    $splitPoint.Collapse(0)       # collapse to right after "33% "

    $splitPoint.InsertAfter("Mozilla ")
    $splitPoint.Bold = 1
    $splitPoint.Bold = 0
}
